$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Weight" column (H) ---------------------------------------------

# Header cell H2: copy the formatting of the neighboring header cell (G2)
# so it picks up the same style as the rest of row 2, then set its text.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = "Weight"

# H1 belongs to the title band (row 1) - match the formatting used by the
# rest of that row before extending the merge over it.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data rows: copy formatting from the existing numeric column in the same
# row (D) so each H cell picks up the matching border/number style, then
# fill in the card weight value.
$ws.Range("D3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = 90

$ws.Range("D4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = 50

$ws.Range("D5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = 50

$ws.Range("D6").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H6").Value = 10

$ws.Range("D7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = 70

$ws.Range("D8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value = 40

# Extend the title merge to cover the new column: A1:G1 -> A1:H1
$ws.Range("A1:G1").UnMerge()
$ws.Range("A1:H1").Merge()

$excel.CutCopyMode = $false
